# Update the "Plot Coordinates" column on the Compartments sheet with the
# improved full cascade coordinate layout, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compartments")

$coords = @{
    2  = "(-3,5)"
    3  = "(3,5)"
    4  = "(0,4)"
    5  = "(-2,3)"
    6  = "(-6,3)"
    7  = "(2,3)"
    8  = "(6,3)"
    9  = "(-8,2)"
    10 = "(-10,1)"
    11 = "(-10,-1)"
    12 = "(-8,1)"
    13 = "(-8,-1)"
    14 = "(-6,1)"
    15 = "(-6,-1)"
    16 = "(0,2)"
    17 = "(-2,1)"
    18 = "(-2,-1)"
    19 = "(0,1)"
    20 = "(0,-1)"
    21 = "(2,1)"
    22 = "(2,-1)"
    23 = "(8,2)"
    24 = "(6,1)"
    25 = "(6,-1)"
    26 = "(8,1)"
    27 = "(8,-1)"
    28 = "(10,1)"
    29 = "(10,-1)"
    30 = "(-4,0)"
    31 = "(4,0)"
}

foreach ($row in $coords.Keys) {
    $ws.Cells.Item($row, 3).Value = $coords[$row]
}

$ws.Range("J15").Select()
